# Swap the content of the two record blocks (rows 18-20) and (rows 21-23)
# on the "Artfynd" sheet. Only the columns that actually carry
# record-specific data are touched: A,B,D,E,F,G,H,K,L,M,N,Q,R,AC.
# (P,S,T,U,V,W,Y,AA,AD,AE,AG,AT,AW,AX,AY are identical between the two
# blocks and are left untouched, avoiding Excel's autodetect-as-date
# behaviour on the "2026-01-25" text cells in Y/AA.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","D","E","F","G","H","K","L","M","N","Q","R","AC")

$pairs = @(
    @(18, 21),
    @(19, 22),
    @(20, 23)
)

foreach ($pair in $pairs) {
    $rowTop = $pair[0]
    $rowBot = $pair[1]

    # Snapshot both rows' values for the relevant columns before writing
    # anything, so the swap doesn't clobber data it still needs to read.
    $topVals = @{}
    $botVals = @{}
    foreach ($col in $cols) {
        $topVals[$col] = $ws.Range("$col$rowTop").Value2
        $botVals[$col] = $ws.Range("$col$rowBot").Value2
    }

    foreach ($col in $cols) {
        $newTop = $botVals[$col]
        $newBot = $topVals[$col]

        if ($null -eq $newTop) {
            $ws.Range("$col$rowTop").ClearContents()
        } else {
            $ws.Range("$col$rowTop").Value = $newTop
        }

        if ($null -eq $newBot) {
            $ws.Range("$col$rowBot").ClearContents()
        } else {
            $ws.Range("$col$rowBot").Value = $newBot
        }
    }
}
